$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Wlan router" row: Leistung [W] (D15) bumped from 5W to 10W
$ws.Range("D15").Value = 10

# "Raspberry (5W converter)" row: now actually present -> Anzahl (B16) 0 -> 1,
# Leitungsquerschnitt (G16) 0 -> 4
$ws.Range("B16").Value = 1
$ws.Range("G16").Value = 4

# Move the view: selection now on H22, zoomed out to 80% (and scrolled back to top-left)
$ws.Range("H22").Select()
$excel.ActiveWindow.Zoom = 80
